# Auto-generated edit script applying numeric updates per the commit diff.
# Workbook contains 8 per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# each backed by an Excel Table (Table_<JOB>) with columns:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 30910.795
$ws.Range("J112").Value = 45072.652
$ws.Range("L112").Value = 135217.956
$ws.Range("N112").Value = -137433.956

$ws.Range("H116").Value = 4518.75
$ws.Range("I116").Value = 3600
$ws.Range("K116").Value = 3600
$ws.Range("M116").Value = -158

$ws.Range("H117").Value = 371040000
$ws.Range("J117").Value = 371040000
$ws.Range("L117").Value = 371040000
$ws.Range("N117").Value = -371049178

$ws.Range("H134").Value = 88999
$ws.Range("J134").Value = 88999
$ws.Range("L134").Value = 88999
$ws.Range("N134").Value = -99139

$ws.Range("H141").Value = 7695.4165
$ws.Range("I141").Value = 8557.5
$ws.Range("K141").Value = 25672.5
$ws.Range("M141").Value = -20492.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4857.86
$ws.Range("I32").Value = 5316.2046
$ws.Range("K32").Value = 5316.2046
$ws.Range("M32").Value = -5029.2046

$ws.Range("H45").Value = 7404.364
$ws.Range("I45").Value = 7429.8
$ws.Range("J45").Value = 7383.1665
$ws.Range("K45").Value = 7429.8
$ws.Range("L45").Value = 7383.1665
$ws.Range("M45").Value = -7052.8
$ws.Range("N45").Value = -8137.1665

$ws.Range("H61").Value = 4292.0654
$ws.Range("I61").Value = 3639.9714
$ws.Range("J61").Value = 6366.909
$ws.Range("K61").Value = 3639.9714
$ws.Range("L61").Value = 6366.909
$ws.Range("M61").Value = -3427.9714
$ws.Range("N61").Value = -6790.909

$ws.Range("H97").Value = 19684.773
$ws.Range("I97").Value = 9967.883
$ws.Range("J97").Value = 52722.2
$ws.Range("K97").Value = 9967.883
$ws.Range("L97").Value = 52722.2
$ws.Range("M97").Value = -9471.883
$ws.Range("N97").Value = -53714.2

$ws.Range("H102").Value = 12110.5
$ws.Range("I102").Value = 12969.059
$ws.Range("J102").Value = 10488.777
$ws.Range("K102").Value = 12969.059
$ws.Range("L102").Value = 10488.777
$ws.Range("M102").Value = -11347.059
$ws.Range("N102").Value = -13732.777

$ws.Range("H110").Value = 3491.5833
$ws.Range("I110").Value = 2899.8333
$ws.Range("J110").Value = 4083.3333
$ws.Range("K110").Value = 2899.8333
$ws.Range("L110").Value = 4083.3333
$ws.Range("M110").Value = -854.8332999999998
$ws.Range("N110").Value = -8173.3333

$ws.Range("H122").Value = 456480.3
$ws.Range("I122").Value = 2962.625
$ws.Range("K122").Value = 8887.875
$ws.Range("M122").Value = -6437.875

$ws.Range("H136").Value = 4292.0654
$ws.Range("I136").Value = 3639.9714
$ws.Range("J136").Value = 6366.909
$ws.Range("K136").Value = 10919.9142
$ws.Range("L136").Value = 19100.727
$ws.Range("M136").Value = -8369.914199999999
$ws.Range("N136").Value = -24200.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6628
$ws.Range("I86").Value = 6763.6665
$ws.Range("K86").Value = 6763.6665
$ws.Range("M86").Value = -5640.6665

$ws.Range("H89").Value = 6628
$ws.Range("I89").Value = 6763.6665
$ws.Range("K89").Value = 33818.3325
$ws.Range("M89").Value = -28202.3325

$ws.Range("H107").Value = 4820.2856
$ws.Range("I107").Value = 5068.4
$ws.Range("J107").Value = 4200
$ws.Range("K107").Value = 5068.4
$ws.Range("L107").Value = 4200
$ws.Range("M107").Value = -3148.4
$ws.Range("N107").Value = -8040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 456.44446
$ws.Range("J10").Value = 491.66666
$ws.Range("L10").Value = 491.66666
$ws.Range("N10").Value = -769.66666

$ws.Range("H22").Value = 1027.6154
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1159.909
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1159.909
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -1859.909

$ws.Range("H107").Value = 35718896
$ws.Range("I107").Value = 47624704
$ws.Range("J107").Value = 1471.2858
$ws.Range("K107").Value = 47624704
$ws.Range("L107").Value = 1471.2858
$ws.Range("M107").Value = -47622784
$ws.Range("N107").Value = -5311.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 627135.4399999999
$ws.Range("J5").Value = 2002093.6
$ws.Range("L5").Value = 6006280.800000001
$ws.Range("N5").Value = -6006504.800000001

$ws.Range("H116").Value = 111113700
$ws.Range("J116").Value = 1877.5
$ws.Range("L116").Value = 5632.5
$ws.Range("N116").Value = -12516.5

$ws.Range("H120").Value = 166669170
$ws.Range("I120").Value = 166669170
$ws.Range("K120").Value = 500007510
$ws.Range("M120").Value = -500002672

$ws.Range("H124").Value = 12120
$ws.Range("I124").Value = 3500
$ws.Range("J124").Value = 13844
$ws.Range("K124").Value = 10500
$ws.Range("L124").Value = 41532
$ws.Range("M124").Value = -5590
$ws.Range("N124").Value = -51352

$ws.Range("H132").Value = 57553.11
$ws.Range("I132").Value = 895
$ws.Range("K132").Value = 8055
$ws.Range("M132").Value = -5525

$ws.Range("H134").Value = 1833.3334
$ws.Range("I134").Value = 1833.3334
$ws.Range("K134").Value = 5500.0002
$ws.Range("M134").Value = -430.0002000000004

$ws.Range("H135").Value = 627135.4399999999
$ws.Range("J135").Value = 2002093.6
$ws.Range("L135").Value = 18018842.4
$ws.Range("N135").Value = -18023912.4

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H138").Value = 3271
$ws.Range("I138").Value = 3271
$ws.Range("K138").Value = 9813
$ws.Range("M138").Value = -4673

$ws.Range("H139").Value = 1073044.8
$ws.Range("I139").Value = 1155202
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 3465606
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -3460466
$ws.Range("N139").Value = -25280

$ws.Range("H140").Value = 13593.533
$ws.Range("I140").Value = 15069.462
$ws.Range("K140").Value = 45208.386
$ws.Range("M140").Value = -40028.386

$ws.Range("H141").Value = 1030
$ws.Range("I141").Value = 1030
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3090
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2090
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6011000
$ws.Range("I14").Value = 6011000
$ws.Range("K14").Value = 6011000
$ws.Range("M14").Value = -6010832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7904.048
$ws.Range("J22").Value = 3170.6667
$ws.Range("L22").Value = 3170.6667
$ws.Range("N22").Value = -3760.6667

$ws.Range("H27").Value = 7904.048
$ws.Range("J27").Value = 3170.6667
$ws.Range("L27").Value = 3170.6667
$ws.Range("N27").Value = -3384.6667

$ws.Range("H40").Value = 18818.625
$ws.Range("I40").Value = 20661.217
$ws.Range("J40").Value = 14109.777
$ws.Range("K40").Value = 20661.217
$ws.Range("L40").Value = 14109.777
$ws.Range("M40").Value = -20525.217
$ws.Range("N40").Value = -14381.777

$ws.Range("H50").Value = 40023.332
$ws.Range("I50").Value = 25037.5
$ws.Range("K50").Value = 25037.5
$ws.Range("M50").Value = -24400.5

$ws.Range("H54").Value = 24492
$ws.Range("J54").Value = 24492
$ws.Range("L54").Value = 24492
$ws.Range("N54").Value = -25780

$ws.Range("H122").Value = 3820.2258
$ws.Range("J122").Value = 4513
$ws.Range("L122").Value = 13539
$ws.Range("N122").Value = -18439

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64499.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 64499.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 64499.5
$ws.Range("N46").Value = -64961.5
$ws.Range("M46").ClearContents()

$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774

$ws.Range("H113").Value = 4384.5713
$ws.Range("I113").Value = 2172.75
$ws.Range("K113").Value = 6518.25
$ws.Range("M113").Value = -4348.25

$ws.Range("H134").Value = 64499.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 64499.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 193498.5
$ws.Range("N134").Value = -198568.5
$ws.Range("M134").ClearContents()
